$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. '68.575.85', '1.00') that
# Excel would otherwise auto-convert to a number when assigned via .Value.
# Force those cells to Text format first, then restore the default 'Normal'
# style afterwards so no stray number-format style gets left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '68.575.85'
$ws.Range('E2').Value = '  +0.51%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.705.29'
$ws.Range('E3').Value = '  +2.17%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '599.06'
$ws.Range('E5').Value = '  +0.24%  '

# Row 6
Set-TextValue $ws.Range('D6') '160.65'
$ws.Range('E6').Value = '  +2.58%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
Set-TextValue $ws.Range('D9') '2.703.97'
$ws.Range('E9').Value = '  +2.17%  '

# Row 10
$ws.Range('E10').Value = '  +0.15%  '

# Row 11
$ws.Range('E11').Value = '  -0.30%  '

# Row 12
Set-TextValue $ws.Range('D12') '5.31'
$ws.Range('E12').Value = '  +1.02%  '

# Row 13
Set-TextValue $ws.Range('D13') '0.361'
$ws.Range('E13').Value = '  +2.80%  '

# Row 14
Set-TextValue $ws.Range('D14') '28.33'
$ws.Range('E14').Value = '  +1.15%  '

# Row 15
Set-TextValue $ws.Range('D15') '3.193.99'
$ws.Range('E15').Value = '  +2.08%  '

# Row 16
$ws.Range('E16').Value = '  -0.96%  '

# Row 17
Set-TextValue $ws.Range('D17') '68.566.15'
$ws.Range('E17').Value = '  +0.47%  '

# Row 18
Set-TextValue $ws.Range('D18') '2.720.90'
$ws.Range('E18').Value = '  +2.98%  '

# Row 19
Set-TextValue $ws.Range('D19') '11.86'
$ws.Range('E19').Value = '  +4.33%  '

# Row 20
Set-TextValue $ws.Range('D20') '7.67'
$ws.Range('E20').Value = '  +4.51%  '

# Row 21
Set-TextValue $ws.Range('D21') '364.83'
$ws.Range('E21').Value = '  +0.44%  '

# Row 22
Set-TextValue $ws.Range('D22') '4.55'
$ws.Range('E22').Value = '  +2.81%  '

# Row 23
$ws.Range('E23').Value = '  +2.23%  '

# Row 24
$ws.Range('E24').Value = '  +2.50%  '

# Row 25
Set-TextValue $ws.Range('D25') '74.12'
$ws.Range('E25').Value = '  -1.54%  '

# Row 26
$ws.Range('E26').Value = '  -0.05%  '

# Row 27
Set-TextValue $ws.Range('D27') '9.93'
$ws.Range('E27').Value = '  +1.81%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.836.81'
$ws.Range('E28').Value = '  +1.98%  '

# Row 29
$ws.Range('E29').Value = '  +0.92%  '

# Row 30
Set-TextValue $ws.Range('D30') '591.99'
$ws.Range('E30').Value = '  +6.07%  '

# Row 31
Set-TextValue $ws.Range('D31') '1.00'
$ws.Range('E31').Value = '  +0.08%  '

# Row 32
$ws.Range('E32').Value = '  +2.25%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.45'
$ws.Range('E33').Value = '  +2.63%  '

# Row 35
$ws.Range('E35').Value = '  +3.22%  '

# Row 36
$ws.Range('E36').Value = '  +5.31%  '

# Row 37
$ws.Range('E37').Value = '  -0.02%  '

# Row 38
Set-TextValue $ws.Range('D38') '161.39'
$ws.Range('E38').Value = '  +0.20%  '

# Row 39
Set-TextValue $ws.Range('D39') '19.85'
$ws.Range('E39').Value = '  +0.82%  '

# Row 40
$ws.Range('E40').Value = '  +2.22%  '

# Row 41
$ws.Range('E41').Value = '  +1.97%  '

# Row 42
Set-TextValue $ws.Range('D42') '5.38'
$ws.Range('E42').Value = '  +1.22%  '

# Row 43
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D43') '18.00'
$ws.Range('E43').Value = '  +1.16%  '

# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D44') '2.67'
$ws.Range('E44').Value = '  +2.70%  '

# Row 45
$ws.Range('E45').Value = '  +0.05%  '

# Row 46
$ws.Range('E46').Value = '  -5.61%  '

# Row 47
Set-TextValue $ws.Range('D47') '157.99'
$ws.Range('E47').Value = '  -0.58%  '

# Row 48
$ws.Range('E48').Value = '  +5.64%  '

# Row 49
$ws.Range('E49').Value = '  +5.17%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.605'
$ws.Range('E50').Value = '  +7.53%  '

# Row 51
Set-TextValue $ws.Range('D51') '22.12'
$ws.Range('E51').Value = '  +0.31%  '
